# "Updates to execute RAD Extension Payment Type."
#
# Only the "Extension Payments" test row (row 4) should be flagged to
# Execute ("Y") going forward, and it gets a refreshed run timestamp.
# The other RAD Payment Type rows no longer carry the "Execute" flag.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the execution timestamp for the Extension Payments row (row 4).
$ws.Range("B4").Value = "Wed Mar 20 23:05:30 EDT 2024"

# Remove the "Execute" (Y) flag from every row except row 4 (Extension
# Payments), which keeps it.
$ws.Range("C2").Clear()
$ws.Range("C3").Clear()
$ws.Range("C5").Clear()
$ws.Range("C6").Clear()
$ws.Range("C7").Clear()
